$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump Version / Date, and insert a new "Jurisdiction"
#    property row right after "Contact" (pushing Description/Purpose/
#    Copyright/Immutable down by one row).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Snapshot the rows that need to shift down before we overwrite anything.
$a11 = $ws.Range("A11").Value2
$b11 = $ws.Range("B11").Value2
$a12 = $ws.Range("A12").Value2
$b12 = $ws.Range("B12").Value2
$a13 = $ws.Range("A13").Value2
$b13 = $ws.Range("B13").Value2
$a14 = $ws.Range("A14").Value2
$b14 = $ws.Range("B14").Value2

# Give the brand-new row 15 the same formatting as the existing data rows
# before writing values into it (otherwise it would pick up the default style).
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A15").Value = $a14
$ws.Range("B15").Value = $b14
$ws.Range("A14").Value = $a13
$ws.Range("B14").Value = $b13
$ws.Range("A13").Value = $a12
$ws.Range("B13").Value = $b12
$ws.Range("A12").Value = $a11
$ws.Range("B12").Value = $b11
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# ---------------------------------------------------------------------------
# 2. Rename the second sheet.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Include from RoleClass")
$ws2.Name = "Include #0"
